$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vm_pu values for rows 2-25 (bus voltage results for the 380 kV case),
# columns B,C,D,E,F,I,J,K,L,M,N (G stays 1, H stays blank).
$columns = @("B","C","D","E","F","I","J","K","L","M","N")

$data = @{
    2 = @(1.02, 1.011247580710721, 1.028009207662965, 1.013424178737686, 1.0251154594127, 1.029811799039262, 1.016497153369578, 1.030827272708938, 1.016285372698983, 1.02794197529134, 1.017940695921216)
    3 = @(1.02, 1.012462214098829, 1.028385681852181, 1.014462673929747, 1.026451009472275, 1.029860433180001, 1.017343547348405, 1.031012757685353, 1.017127943177975, 1.029083325371322, 1.018788291876576)
    4 = @(1.02, 1.013248101984948, 1.028629008135119, 1.01513491281257, 1.027314502582981, 1.029890287381812, 1.017890722470564, 1.031131757354537, 1.017672814189411, 1.029820632288004, 1.01933624405017)
    5 = @(1.02, 1.013578477192167, 1.028731234194849, 1.015417586367324, 1.027677350816842, 1.02990245080704, 1.018120637062317, 1.031181539131005, 1.017901800633994, 1.030130305448863, 1.019566485147015)
    6 = @(1.02, 1.013633948010818, 1.028748394340917, 1.015465052312513, 1.02773826507482, 1.029904470372765, 1.018159233867412, 1.031189883270182, 1.017940243968708, 1.030182283951165, 1.019605136764001)
    7 = @(1.02, 1.013252516518816, 1.028630374355475, 1.015138689656, 1.027319351621112, 1.029890451432565, 1.017893795061747, 1.031132423507885, 1.017675874219735, 1.029824771295508, 1.019339321004785)
    8 = @(1.02, 1.011658085960854, 1.028136495279595, 1.01377508939948, 1.025566960719816, 1.029828569302794, 1.016783300264107, 1.03089016905827, 1.016570192123274, 1.028327954371427, 1.018227249177157)
    9 = @(1.02, 1.008847911892502, 1.027264170815833, 1.01137420513929, 1.022473567584923, 1.02970717824176, 1.014822582809363, 1.030455511838382, 1.014619265125886, 1.02568091103078, 1.01626374727873)
    10 = @(1.02, 1.006973906831952, 1.026681353410436, 1.009774841449518, 1.020407469441405, 1.029617982363476, 1.013512731437356, 1.030160580180597, 1.013316832849115, 1.023909719252517, 1.014952035767532)
    11 = @(1.02, 1.006162273272842, 1.026428709113369, 1.009082569826139, 1.019511878098101, 1.029577405125989, 1.012944887639038, 1.030031661521243, 1.01275241572005, 1.023141201772315, 1.014383385565888)
    12 = @(1.02, 1.005860766583414, 1.026334825467372, 1.00882546706854, 1.0191790683791, 1.02956203975037, 1.012733862629327, 1.029983594512264, 1.012542695918734, 1.022855500021797, 1.014172060876458)
    13 = @(1.02, 1.005925442194816, 1.026354965622651, 1.008880614804119, 1.019250463910383, 1.029565348934628, 1.01277913288184, 1.029993913216344, 1.012587684730376, 1.022916794929088, 1.014217395417922)
    14 = @(1.02, 1.00613735126358, 1.026420949473536, 1.009061316872573, 1.019484370965883, 1.029576140997166, 1.012927446349013, 1.030027691970798, 1.012735081656939, 1.023117590491535, 1.014365919507231)
    15 = @(1.02, 1.006267911382106, 1.026461599015739, 1.009172658327858, 1.019628469212649, 1.029582751508626, 1.013018813525018, 1.030048480240978, 1.012825888420619, 1.023241275359705, 1.014457416435105)
    16 = @(1.02, 1.007027768271366, 1.026698114827572, 1.009820790631153, 1.020466886431473, 1.029620634198533, 1.013550403087732, 1.030169110661264, 1.013354281613189, 1.023960689704411, 1.014989760915973)
    17 = @(1.02, 1.007504356738725, 1.026846401148775, 1.010227416135606, 1.020992544728224, 1.02964387400521, 1.013883674552005, 1.030244455419853, 1.013685605566608, 1.024411534399127, 1.015323505663952)
    18 = @(1.02, 1.007782325871752, 1.026932866863317, 1.010464619452835, 1.02129906019447, 1.029657240744564, 1.014078001693758, 1.030288285796009, 1.013878817406158, 1.024674352249912, 1.015518108772535)
    19 = @(1.02, 1.007877103401476, 1.026962344769694, 1.010545504040892, 1.021403558523707, 1.029661766443858, 1.0141442513515, 1.030303210951503, 1.013944690270353, 1.024763940607988, 1.015584452512391)
    20 = @(1.02, 1.007453225074932, 1.02683049422893, 1.010183786439702, 1.020936156122776, 1.029641400102964, 1.013847924360864, 1.030236383728821, 1.013650062187632, 1.024363178797766, 1.015287704703439)
    21 = @(1.02, 1.006074950174649, 1.026401519971977, 1.009008103590018, 1.019415495259401, 1.029572971092956, 1.012883774575872, 1.030017749953829, 1.012691678886274, 1.0230584678547, 1.014322185715162)
    22 = @(1.02, 1.00520819973406, 1.026131574537345, 1.008269122599374, 1.018458541962556, 1.02952825104067, 1.012276981102651, 1.029879240253532, 1.012088698541257, 1.022236753526111, 1.013714530525609)
    23 = @(1.02, 1.005667698154082, 1.026274699128386, 1.008660850402709, 1.018965923052155, 1.029552118580874, 1.012598710785068, 1.029952765645281, 1.012408388965879, 1.02267249260813, 1.014036717101078)
    24 = @(1.02, 1.00747632932486, 1.026837681965616, 1.010203500742382, 1.020961635990909, 1.029642518535893, 1.013864078532558, 1.030240031340096, 1.013666122846035, 1.024385029082883, 1.01530388181591)
    25 = @(1.02, 1.009574495204962, 1.027489918915731, 1.011994668757404, 1.023273947664625, 1.029740020505447, 1.015329945290461, 1.030568795504558, 1.015123940135079, 1.026366370212686, 1.016771830272745)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range("$($columns[$i])$row").Value = $values[$i]
    }
}

Write-Output "Updated vm_pu results for rows 2-25 (380 kV case)"